# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet to the latest snapshot, and fixes the ordering for two coin pairs
# (PancakeSwap/Dai swap back to Dai/PancakeSwap at rows 25-26, and
# Monero/Maker swap to Maker/Monero at rows 42-43) whose rank changed.
#
# Column D holds prices formatted as plain text (e.g. "43.082.19",
# "1.00", "0.999") rather than numbers, so assigning a bare numeric-
# looking string via .Value would make Excel auto-coerce it into a
# real number (dropping the original text formatting, e.g. "1.00" ->
# 1). To avoid that we prefix the literal with a leading apostrophe
# (forces text entry, like typing '1.00 into a cell) and then reset
# the cell's Style back to Normal so no stray number-format/quote-
# prefix formatting is left on the cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.095.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = "'2.355.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.52%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'302.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = "'95.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = "'34.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("D11").Value = "'0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = "'18.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.87%  '
$ws.Range("E13").Value = '  +3.26%  '
$ws.Range("D14").Value = "'6.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("D15").Value = "'2.720.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").Value = "'2.353.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").Value = "'43.070.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("E20").Value = '  +4.30%  '
$ws.Range("D21").Value = "'0.0₃0888"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = "'68.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").Value = "'235.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = "'2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = "'2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").Value = "'31.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("D32").Value = "'5.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").Value = "'0.0725"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").Value = "'17.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.66%  '
$ws.Range("E35").Value = '  -1.97%  '
$ws.Range("E36").Value = '  +4.90%  '
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("E39").Value = '  +9.93%  '
$ws.Range("E40").Value = '  +2.24%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = "'1.943.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = "'102.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -37.54%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = "'2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.76%  '
$ws.Range("D46").Value = "'9.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.58%  '
$ws.Range("E47").Value = '  -0.75%  '
$ws.Range("D48").Value = "'2.583.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = "'52.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("E50").Value = '  -4.21%  '
$ws.Range("D51").Value = "'72.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.43%  '
